# Auto-generated edit script: updates numeric price/profit data cells
# across several worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# a refreshed market-price snapshot. A couple of now-stale cells (N111/N112
# on GSM) are cleared entirely rather than zeroed.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H21" = 17777.555
    "I21" = 8000
    "J21" = 18999.75
    "K21" = 8000
    "L21" = 18999.75
    "M21" = -7532
    "N21" = -19935.75
    "H23" = 17777.555
    "I23" = 8000
    "J23" = 18999.75
    "K23" = 8000
    "L23" = 18999.75
    "M23" = -7766
    "N23" = -19467.75
    "H29" = 29.8
    "I29" = 29.8
    "K29" = 89.40000000000001
    "M29" = 191.6
    "H38" = 786.0417
    "I38" = 200.25
    "J38" = 1371.8334
    "K38" = 600.75
    "L38" = 4115.5002
    "M38" = -228.75
    "N38" = -4859.5002
    "H58" = 2050.2
    "I58" = 1475.3
    "J58" = 3200
    "K58" = 4425.9
    "L58" = 9600
    "M58" = -4275.9
    "N58" = -9900
    "H87" = 25858.5
    "J87" = 25858.5
    "L87" = 25858.5
    "N87" = -28354.5
    "H90" = 25858.5
    "J90" = 25858.5
    "L90" = 77575.5
    "N90" = -90055.5
    "H112" = 5078.049
    "J112" = 5582.162
    "L112" = 16746.486
    "N112" = -18962.486
    "H132" = 2438.8086
    "I132" = 2561.5676
    "J132" = 1984.6
    "K132" = 7684.702799999999
    "L132" = 5953.799999999999
    "M132" = -5154.702799999999
    "N132" = -11013.8
    "H134" = 79520
    "J134" = 79520
    "L134" = 79520
    "N134" = -89660
    "H141" = 3397.8262
    "I141" = 1503.6111
    "J141" = 10217
    "K141" = 4510.8333
    "L141" = 30651
    "M141" = 669.1666999999998
    "N141" = -41011
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H2" = 1715
    "I2" = 1672.5
    "K2" = 1672.5
    "M2" = -1559.5
    "H74" = 19234010
    "I74" = 2308.5386
    "K74" = 2308.5386
    "M74" = -1434.5386
    "H77" = 19234010
    "I77" = 2308.5386
    "K77" = 11542.693
    "M77" = -7174.692999999999
    "H116" = 1715
    "I116" = 1672.5
    "K116" = 1672.5
    "M116" = 621.5
    "H132" = 1427715.2
    "I132" = 2483.4443
    "J132" = 8553874
    "K132" = 7450.3329
    "L132" = 25661622
    "M132" = -4920.3329
    "N132" = -25666682
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H3" = 1715
    "I3" = 1672.5
    "K3" = 1672.5
    "M3" = -1558.5
    "H107" = 112466.555
    "I107" = 200899.8
    "J107" = 1925
    "K107" = 200899.8
    "L107" = 1925
    "M107" = -198979.8
    "N107" = -5765
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H132" = 19049790
    "I132" = 20835220
    "J132" = 15154308
    "K132" = 62505660
    "L132" = 45462924
    "M132" = -62503130
    "N132" = -45467984
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H34" = 11628472
    "I34" = 200
    "J34" = 12195705
    "K34" = 600
    "L34" = 36587115
    "M34" = -516
    "N34" = -36587283
    "H39" = 2354.4546
    "J39" = 2354.4546
    "L39" = 7063.3638
    "N39" = -7651.3638
    "H55" = 1477.8334
    "I55" = 1300
    "J55" = 1513.4
    "K55" = 3900
    "L55" = 4540.200000000001
    "M55" = -3723
    "N55" = -4894.200000000001
    "H136" = 2972.2307
    "I136" = 1183.9
    "J136" = 8933.333000000001
    "K136" = 3551.7
    "L136" = 26799.999
    "M136" = 1548.3
    "N136" = -36999.999
    "H137" = 35376.65
    "I137" = 8093.2354
    "J137" = 62660.06
    "K137" = 24279.7062
    "L137" = 187980.18
    "M137" = -19179.7062
    "N137" = -198180.18
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H110" = 41111
    "J110" = 41111
    "L110" = 41111
    "N110" = -49291
    "H111" = 0
    "J111" = 0
    "L111" = 0
    "H112" = 0
    "J112" = 0
    "L112" = 0
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}
$ws.Range("N111").ClearContents() | Out-Null
$ws.Range("N112").ClearContents() | Out-Null

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H7" = 5252.5386
    "I7" = 5262.091
    "J7" = 5200
    "K7" = 5262.091
    "L7" = 5200
    "M7" = -5150.091
    "N7" = -5424
    "H126" = 5252.5386
    "I126" = 5262.091
    "J126" = 5200
    "K126" = 15786.273
    "L126" = 15600
    "M126" = -13316.273
    "N126" = -20540
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H54" = 9921.5
    "J54" = 9921.5
    "L54" = 9921.5
    "N54" = -10961.5
    "H81" = 3279.1875
    "J81" = 2299.75
    "L81" = 4599.5
    "N81" = -6721.5
    "H84" = 3279.1875
    "J84" = 2299.75
    "L84" = 22997.5
    "N84" = -33605.5
}
foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}
